$d = $word.ActiveDocument

# --- Locate the two paragraphs affected by this edit -----------------------
# Paragraph A: "Prompt: - Based on SPEC.md sections ... ( ... Typography ...)."
#              The leading "Prompt: " stays un-struck; everything after it in
#              this paragraph becomes struck-through.
# Paragraph B: "Acceptance criteria: - Each property maps per SPEC.md ..."
#              The whole paragraph (including its paragraph mark) becomes
#              struck-through.

$promptPara = $null
$acceptancePara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($promptPara -eq $null -and $t.StartsWith("Prompt: - Based on SPEC.md sections")) {
        $promptPara = $p
        $acceptancePara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($promptPara -eq $null) {
    throw "Could not locate the target 'Prompt: - Based on SPEC.md ...' paragraph"
}

# --- Paragraph A: strike everything after the literal "Prompt: " prefix ----
$paraRange = $promptPara.Range
$prefixLen = "Prompt: ".Length
$afterPrefixStart = $paraRange.Start + $prefixLen

# Stop just before the paragraph mark (Range.End of a Paragraphs item points
# past the trailing pilcrow), so this only covers the visible run text.
$strikeRange = $d.Range($afterPrefixStart, $paraRange.End - 1)
$strikeRange.Font.StrikeThrough = 1

# --- Paragraph B: strike the whole paragraph, including the mark -----------
$acceptancePara.Range.Font.StrikeThrough = 1

Write-Output "strikethrough applied"
